# Collapse the 4-column API-key-management sheet down to a single
# "input_KeyName" column: drop columns A-C (their content + widths),
# which shifts former column D (header "input_KeyName" + blank row 2)
# into column A, preserving D1's header style and D2's blank value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:C").EntireColumn.Delete()
